$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.723.44"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "'2.207.09"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'229.81"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'60.49"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  -3.16%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'2.532.93"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "'15.38"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.794"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'5.57"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "'2.215.40"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "'41.686.33"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'0.0₃0903"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'71.98"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'6.04"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").Value = "'241.78"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  -4.75%  "
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("D28").Value = "'168.34"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("D31").Value = "'19.72"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("E32").Value = "  -8.08%  "
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "'0.0644"
$ws.Range("E36").Value = "  +4.24%  "
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").Value = "'6.27"
$ws.Range("E38").Value = "  -5.87%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").Value = "'0.000244"
$ws.Range("E40").Value = "  +6.75%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").Value = "'8.69"
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("D46").Value = "'4.38"
$ws.Range("E46").Value = "  -9.58%  "
$ws.Range("D47").Value = "'96.41"
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").Value = "'1.461.64"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("E50").Value = "  -3.87%  "
$ws.Range("E51").Value = "  +0.08%  "
